$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.175.08'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.27%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.310.67'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.80%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.37'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.72%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.09'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.99%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.517'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +4.56%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.59'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +9.15%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0794'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.35%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.73'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.93%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.99'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.49%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.670.30'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.66%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.259.67'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.83%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.808'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.80%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.054.61'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.42'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +7.67%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.26'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.08%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.43%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.27'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.15%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.39%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.26'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +14.73%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.31%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.22%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.14'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.72%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.49%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '35.11'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.96%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '169.54'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.71%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.22'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.75%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.13%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.06'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.05%  '

$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.94'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.16%  '

$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.73'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.99%  '

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.32%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.57%  '

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.66%  '

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.64%  '

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.50%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.32%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.996.83'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.01%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0290'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.76%  '

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.47%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.29'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +5.13%  '

$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.93'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.92%  '

$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.76'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.37%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '56.25'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +6.62%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.56'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.97%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.536.70'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.68%  '

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.28%  '
